$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# All edits are applied by locating the exact target paragraph (by its
# current text, and - where the text is ambiguous - its style too, so the
# two paragraphs that happen to share the same old text, the H1 title and
# the later bold "Play Joker..." line, are never confused with one
# another) and then rewriting that single paragraph's Range via InsertXML.
# InsertXML replaces only the content of the given Range, which keeps
# every other paragraph/run completely untouched and preserves
# formatting-carrying empty <w:r/> runs that a plain Find/Replace would
# otherwise silently merge away.

# 1. Main page H1 title.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Play Joker Super Reels Free Slot Game") -and $p.Style.NameLocal -eq "Heading 1") {
        $p.Range.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""Heading1""/></w:pPr><w:r><w:t>Play Joker Super Reels for Free!</w:t></w:r></w:p>")
        break
    }
}

# 2-7. "What we like" / "What we don't like" bullet paragraphs.
function Set-BulletParagraphText($needle, $newText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            $xml = "<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>$newText</w:t></w:r></w:p>"
            $p.Range.InsertXML($xml)
            break
        }
    }
}

Set-BulletParagraphText "Two gameplay grids increase variety" "Two game grids add excitement"
Set-BulletParagraphText "Autoplay feature saves time and effort" "Super Game offers additional symbols and prizes"
Set-BulletParagraphText "Super Game with up to 64 paylines offers exciting rewards" "High volatility and fair RTP"
Set-BulletParagraphText "Suitable for all types of players with variable bets" "Suitable for all types of players"
Set-BulletParagraphText "High volatility means wins may be infrequent" "Wins may be infrequent due to high volatility"
Set-BulletParagraphText "Limited to five fixed paylines in the basic game" "Limited number of paylines"

# 8. Bold "Play Joker..." paragraph near the end (empty run + bold run).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Play Joker Super Reels Free Slot Game") -and $p.Style.NameLocal -ne "Heading 1") {
        $p.Range.InsertXML("<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Joker Super Reels for Free!</w:t></w:r></w:p>")
        break
    }
}

# 9. Italic meta-description paragraph at the very end (empty run + italic run).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Review of Joker Super Reels")) {
        $xml = "<w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the exciting gameplay and bonuses of Joker Super Reels. Play for free now!</w:t></w:r></w:p>"
        $p.Range.InsertXML($xml)
        break
    }
}
